$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D entirely (was part of a 4x4 block, now reduced to 3x3)
$ws.Columns("D").Delete()

# Delete row 4 entirely (was part of a 4x4 block, now reduced to 3x3)
$ws.Rows("4").Delete()

# Update remaining 3x3 block values to the new pattern (all 1s except diagonal 0s)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 1

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0

# Update selection to match the diff (activeCell D7)
$ws.Range("D7").Select()
